$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 67, shifting existing rows 67..171 down to 68..172
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row 67 with the new data record
$ws.Cells.Item(67, 1).Value = 7
$ws.Cells.Item(67, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(67, 3).Value = "Ñuble"
$ws.Cells.Item(67, 4).Value = 44482
$ws.Cells.Item(67, 5).Value = 16
$ws.Cells.Item(67, 6).Value = 100112023
$ws.Cells.Item(67, 7).Value = "Brócoli"
$ws.Cells.Item(67, 8).Value = "Sin especificar"
$ws.Cells.Item(67, 9).Value = "Primera"
$ws.Cells.Item(67, 10).Value = 300
$ws.Cells.Item(67, 11).Value = 750
$ws.Cells.Item(67, 12).Value = 800
$ws.Cells.Item(67, 13).Value = 775
$ws.Cells.Item(67, 14).Value = "`$/unidad"
$ws.Cells.Item(67, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(67, 16).Value = 775
$ws.Cells.Item(67, 17).Value = 1
$ws.Cells.Item(67, 18).Value = "Hortaliza"
